$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2607.7
$ws.Cells.Item(32, 9).Value = 1597.3334
$ws.Cells.Item(32, 10).Value = 3040.7144
$ws.Cells.Item(32, 11).Value = 1597.3334
$ws.Cells.Item(32, 12).Value = 3040.7144
$ws.Cells.Item(32, 13).Value = -1271.3334
$ws.Cells.Item(32, 14).Value = -3692.7144
$ws.Cells.Item(101, 8).Value = 20000466
$ws.Cells.Item(101, 9).Value = 25000434
$ws.Cells.Item(101, 10).Value = 600
$ws.Cells.Item(101, 11).Value = 75001302
$ws.Cells.Item(101, 12).Value = 1800
$ws.Cells.Item(101, 13).Value = -74999680
$ws.Cells.Item(101, 14).Value = -5044
$ws.Cells.Item(112, 8).Value = 2970
$ws.Cells.Item(112, 10).Value = 3067.1
$ws.Cells.Item(112, 12).Value = 9201.299999999999
$ws.Cells.Item(112, 14).Value = -11417.3
$ws.Cells.Item(137, 8).Value = 1228.2
$ws.Cells.Item(137, 9).Value = 1228.2
$ws.Cells.Item(137, 11).Value = 3684.6
$ws.Cells.Item(137, 13).Value = -1134.6
$ws.Cells.Item(138, 8).Value = 2796.8125
$ws.Cells.Item(138, 10).Value = 2577.2
$ws.Cells.Item(138, 12).Value = 7731.599999999999
$ws.Cells.Item(138, 14).Value = -18011.6

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1265.6666
$ws.Cells.Item(2, 9).Value = 700
$ws.Cells.Item(2, 11).Value = 700
$ws.Cells.Item(2, 13).Value = -587
$ws.Cells.Item(21, 8).Value = 1699
$ws.Cells.Item(21, 9).Value = 1699
$ws.Cells.Item(21, 11).Value = 1699
$ws.Cells.Item(21, 13).Value = -1325
$ws.Cells.Item(110, 8).Value = 7400639.5
$ws.Cells.Item(110, 9).Value = 9250674
$ws.Cells.Item(110, 11).Value = 9250674
$ws.Cells.Item(110, 13).Value = -9248629
$ws.Cells.Item(116, 8).Value = 1265.6666
$ws.Cells.Item(116, 9).Value = 700
$ws.Cells.Item(116, 11).Value = 700
$ws.Cells.Item(116, 13).Value = 1594

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1265.6666
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 11).Value = 700
$ws.Cells.Item(3, 13).Value = -586
$ws.Cells.Item(13, 8).Value = 80000
$ws.Cells.Item(13, 10).Value = 80000
$ws.Cells.Item(13, 12).Value = 80000
$ws.Cells.Item(13, 14).Value = -80336
$ws.Cells.Item(22, 8).Value = 340.76923
$ws.Cells.Item(22, 9).Value = 348.18182
$ws.Cells.Item(22, 11).Value = 348.18182
$ws.Cells.Item(22, 13).Value = -175.18182
$ws.Cells.Item(50, 8).Value = 88000
$ws.Cells.Item(50, 10).Value = 88000
$ws.Cells.Item(50, 12).Value = 88000
$ws.Cells.Item(50, 14).Value = -89148
$ws.Cells.Item(87, 8).Value = 126832.664
$ws.Cells.Item(87, 10).Value = 126832.664
$ws.Cells.Item(87, 12).Value = 126832.664
$ws.Cells.Item(87, 14).Value = -129328.664
$ws.Cells.Item(90, 8).Value = 126832.664
$ws.Cells.Item(90, 10).Value = 126832.664
$ws.Cells.Item(90, 12).Value = 380497.992
$ws.Cells.Item(90, 14).Value = -392977.992
$ws.Cells.Item(99, 8).Value = 1625.6666
$ws.Cells.Item(99, 9).Value = 1233
$ws.Cells.Item(99, 11).Value = 1233
$ws.Cells.Item(99, 13).Value = 265

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 5250
$ws.Cells.Item(3, 9).Value = 5250
$ws.Cells.Item(3, 11).Value = 5250
$ws.Cells.Item(3, 13).Value = -5137
$ws.Cells.Item(7, 8).Value = 247.85715
$ws.Cells.Item(7, 9).Value = 61.666668
$ws.Cells.Item(7, 11).Value = 61.666668
$ws.Cells.Item(7, 13).Value = 51.333332
$ws.Cells.Item(19, 8).Value = 33.2
$ws.Cells.Item(19, 9).Value = 33.2
$ws.Cells.Item(19, 11).Value = 33.2
$ws.Cells.Item(19, 13).Value = 136.8
$ws.Cells.Item(22, 8).Value = 641
$ws.Cells.Item(22, 9).Value = 465.2857
$ws.Cells.Item(22, 11).Value = 465.2857
$ws.Cells.Item(22, 13).Value = -115.2857
$ws.Cells.Item(24, 8).Value = 33.2
$ws.Cells.Item(24, 9).Value = 33.2
$ws.Cells.Item(24, 11).Value = 33.2
$ws.Cells.Item(24, 13).Value = 136.8
$ws.Cells.Item(58, 8).Value = 1588
$ws.Cells.Item(58, 9).Value = 1469.1666
$ws.Cells.Item(58, 11).Value = 1469.1666
$ws.Cells.Item(58, 13).Value = -1266.1666
$ws.Cells.Item(105, 8).Value = 2886.56
$ws.Cells.Item(105, 9).Value = 2319.923
$ws.Cells.Item(105, 11).Value = 2319.923
$ws.Cells.Item(105, 13).Value = -572.9229999999998
$ws.Cells.Item(122, 8).Value = 1507.3636
$ws.Cells.Item(122, 9).Value = 1507.3636
$ws.Cells.Item(122, 11).Value = 4522.0908
$ws.Cells.Item(122, 13).Value = -2072.0908
$ws.Cells.Item(136, 8).Value = 1588
$ws.Cells.Item(136, 9).Value = 1469.1666
$ws.Cells.Item(136, 11).Value = 4407.4998
$ws.Cells.Item(136, 13).Value = -1857.4998

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 838.25
$ws.Cells.Item(7, 9).Value = 800.3333
$ws.Cells.Item(7, 11).Value = 2400.9999
$ws.Cells.Item(7, 13).Value = -2288.9999
$ws.Cells.Item(98, 8).Value = 1741.1428
$ws.Cells.Item(98, 9).Value = 194
$ws.Cells.Item(98, 11).Value = 582
$ws.Cells.Item(98, 13).Value = 916
$ws.Cells.Item(107, 8).Value = 468.5
$ws.Cells.Item(107, 9).Value = 289
$ws.Cells.Item(107, 10).Value = 648
$ws.Cells.Item(107, 11).Value = 867
$ws.Cells.Item(107, 12).Value = 1944
$ws.Cells.Item(107, 13).Value = 1053
$ws.Cells.Item(107, 14).Value = -5784
$ws.Cells.Item(109, 8).Value = 1717.7142
$ws.Cells.Item(109, 9).Value = 1804
$ws.Cells.Item(109, 10).Value = 1200
$ws.Cells.Item(109, 11).Value = 5412
$ws.Cells.Item(109, 12).Value = 3600
$ws.Cells.Item(109, 13).Value = -4372
$ws.Cells.Item(109, 14).Value = -5680
$ws.Cells.Item(113, 8).Value = 2066.1428
$ws.Cells.Item(113, 9).Value = 1732.25
$ws.Cells.Item(113, 11).Value = 5196.75
$ws.Cells.Item(113, 13).Value = -3026.75
$ws.Cells.Item(132, 8).Value = 4816.25
$ws.Cells.Item(132, 9).Value = 4015.5
$ws.Cells.Item(132, 10).Value = 5617
$ws.Cells.Item(132, 11).Value = 36139.5
$ws.Cells.Item(132, 12).Value = 50553
$ws.Cells.Item(132, 13).Value = -33609.5
$ws.Cells.Item(132, 14).Value = -55613

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 149.8
$ws.Cells.Item(2, 9).Value = 62.25
$ws.Cells.Item(2, 11).Value = 62.25
$ws.Cells.Item(2, 13).Value = 50.75
$ws.Cells.Item(33, 8).Value = 14053400
$ws.Cells.Item(33, 9).Value = 10000000
$ws.Cells.Item(33, 10).Value = 15066750
$ws.Cells.Item(33, 11).Value = 10000000
$ws.Cells.Item(33, 12).Value = 15066750
$ws.Cells.Item(33, 13).Value = -9999748
$ws.Cells.Item(33, 14).Value = -15067254
$ws.Cells.Item(70, 8).Value = 6460.7144
$ws.Cells.Item(70, 9).Value = 5503.5
$ws.Cells.Item(70, 11).Value = 5503.5
$ws.Cells.Item(70, 13).Value = -5233.5
$ws.Cells.Item(73, 8).Value = 6460.7144
$ws.Cells.Item(73, 9).Value = 5503.5
$ws.Cells.Item(73, 11).Value = 5503.5
$ws.Cells.Item(73, 13).Value = -4567.5
$ws.Cells.Item(80, 8).Value = 1779.6
$ws.Cells.Item(80, 10).Value = 2100
$ws.Cells.Item(80, 12).Value = 2100
$ws.Cells.Item(80, 14).Value = -4096
$ws.Cells.Item(83, 8).Value = 1779.6
$ws.Cells.Item(83, 10).Value = 2100
$ws.Cells.Item(83, 12).Value = 10500
$ws.Cells.Item(83, 14).Value = -20484
$ws.Cells.Item(135, 8).Value = 295000
$ws.Cells.Item(135, 10).Value = 295000
$ws.Cells.Item(135, 12).Value = 295000
$ws.Cells.Item(135, 14).Value = -305140

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 141666.67
$ws.Cells.Item(2, 9).Value = 141666.67
$ws.Cells.Item(2, 11).Value = 141666.67
$ws.Cells.Item(2, 13).Value = -141554.67
$ws.Cells.Item(61, 8).Value = 5333
$ws.Cells.Item(61, 9).Value = 4999.5
$ws.Cells.Item(61, 11).Value = 4999.5
$ws.Cells.Item(61, 13).Value = -4797.5
$ws.Cells.Item(93, 8).Value = 1619.2222
$ws.Cells.Item(93, 9).Value = 1391.7142
$ws.Cells.Item(93, 10).Value = 1864.2307
$ws.Cells.Item(93, 11).Value = 1391.7142
$ws.Cells.Item(93, 12).Value = 1864.2307
$ws.Cells.Item(93, 13).Value = -143.7141999999999
$ws.Cells.Item(93, 14).Value = -4360.2307
$ws.Cells.Item(106, 8).Value = 45226.832
$ws.Cells.Item(106, 10).Value = 45226.832
$ws.Cells.Item(106, 12).Value = 45226.832
$ws.Cells.Item(106, 14).Value = -47750.832
$ws.Cells.Item(113, 8).Value = 5333
$ws.Cells.Item(113, 9).Value = 4999.5
$ws.Cells.Item(113, 11).Value = 4999.5
$ws.Cells.Item(113, 13).Value = -2829.5
$ws.Cells.Item(122, 8).Value = 7379
$ws.Cells.Item(122, 9).Value = 9410.1
$ws.Cells.Item(122, 10).Value = 5347.9
$ws.Cells.Item(122, 11).Value = 28230.3
$ws.Cells.Item(122, 12).Value = 16043.7
$ws.Cells.Item(122, 13).Value = -25780.3
$ws.Cells.Item(122, 14).Value = -20943.7

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 5000
$ws.Cells.Item(3, 10).Value = 5000
$ws.Cells.Item(3, 12).Value = 5000
$ws.Cells.Item(3, 14).Value = -5228
$ws.Cells.Item(7, 8).Value = 18500
$ws.Cells.Item(7, 9).Value = 2000
$ws.Cells.Item(7, 10).Value = 35000
$ws.Cells.Item(7, 11).Value = 2000
$ws.Cells.Item(7, 12).Value = 35000
$ws.Cells.Item(7, 13).Value = -1887
$ws.Cells.Item(7, 14).Value = -35226
$ws.Cells.Item(100, 8).Value = 8334653
$ws.Cells.Item(100, 9).Value = 10001364
$ws.Cells.Item(100, 10).Value = 1099
$ws.Cells.Item(100, 11).Value = 20002728
$ws.Cells.Item(100, 12).Value = 2198
$ws.Cells.Item(100, 13).Value = -20002187
$ws.Cells.Item(100, 14).Value = -3280
$ws.Cells.Item(107, 8).Value = 1500
$ws.Cells.Item(107, 9).Value = 1500
$ws.Cells.Item(107, 11).Value = 4500
$ws.Cells.Item(107, 13).Value = -2580
$ws.Cells.Item(126, 8).Value = 2301.52
$ws.Cells.Item(126, 9).Value = 2377.2104
$ws.Cells.Item(126, 11).Value = 7131.6312
$ws.Cells.Item(126, 13).Value = -4661.6312
